# Apply the "Poisson denoising with rmse=0.373" update to the Experiments sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 53: fill in the metric columns that were previously blank ---
$ws.Range("B53").Value = 0.9239
$ws.Range("C53").Value = 0.9825
$ws.Range("D53").Value = 0.9944
$ws.Range("F53").Value = 0.3552
$ws.Range("G53").Value = 0.084
$ws.Range("I53").Value = 0.036

# --- Row 54 ---
$ws.Range("A54").Value = "DenseDepth (Weighted histogram matching)"
$ws.Range("B54").Value = 0.9239
$ws.Range("C54").Value = 0.9825
$ws.Range("D54").Value = 0.9944
$ws.Range("F54").Value = 0.3552
$ws.Range("G54").Value = 0.084
$ws.Range("I54").Value = 0.036
$ws.Range("L54").Value = "Intensity and Falloff"
$ws.Range("L54").WrapText = $true

# --- Row 55 ---
$ws.Range("A55").Value = "DenseDepth (Weighted histogram matching)"
$ws.Range("B55").Value = 0.9238
$ws.Range("C55").Value = 0.9827
$ws.Range("D55").Value = 0.9947
$ws.Range("F55").Value = 0.3545
$ws.Range("G55").Value = 0.0836
$ws.Range("I55").Value = 0.0358
$ws.Range("L55").Value = "Intensity, Falloff, and DC/Ambient"
$ws.Range("L55").WrapText = $true

# --- Row 56 (write M56 before M55's string so new shared strings are
#     allocated in the same order as the target workbook: "lam=1e-2" (63)
#     then "lam-1e-2" (64)) ---
$ws.Range("M56").Value = "lam=1e-2"
$ws.Range("M55").Value = "lam-1e-2"

$ws.Range("A56").Value = "DenseDepth (Weighted histogram matching)"
$ws.Range("B56").Value = 0.9239
$ws.Range("C56").Value = 0.9828
$ws.Range("D56").Value = 0.9947
$ws.Range("F56").Value = 0.3535
$ws.Range("G56").Value = 0.0836
$ws.Range("I56").Value = 0.0358
$ws.Range("L56").Value = "Intensity, Falloff, DC/Ambient, and Jitter"
$ws.Range("L56").WrapText = $true

# --- Row 57 ---
$ws.Range("A57").Value = "DenseDepth (Weighted histogram matching)"
$ws.Range("B57").Value = 0.9212
$ws.Range("C57").Value = 0.9811
$ws.Range("D57").Value = 0.9936
$ws.Range("F57").Value = 0.3947
$ws.Range("G57").Value = 0.0859
$ws.Range("I57").Value = 0.0366
$ws.Range("L57").Value = "Intensity, Falloff, DC/Ambient, Jitter, and Poisson Noise"
$ws.Range("L57").WrapText = $true
$ws.Range("M57").Value = "lam=1e1"

# --- Row 59 (row 58 intentionally left blank) ---
$ws.Range("A59").Value = "DenseDepth (Weighted histogram matching)"
$ws.Range("B59").Value = 0.9259
$ws.Range("F59").Value = 0.3447
$ws.Range("L59").Value = "Intensity Only (Poissn Denoising) (sid_bins=140) (lam=1e0)"
$ws.Range("L59").WrapText = $true

# --- Row 60 ---
$ws.Range("A60").Value = "DenseDepth (Weighted histogram matching)"
$ws.Range("B60").Value = 0.9206
$ws.Range("C60").Value = 0.9812
$ws.Range("D60").Value = 0.9939
$ws.Range("F60").Value = 0.3726
$ws.Range("G60").Value = 0.0859
$ws.Range("I60").Value = 0.0366
$ws.Range("L60").Value = "Intensity, Falloff, DC/Ambient, Jitter, and Poisson Noise (Poisson Denoising) (sid_bins=140) (lam=1e0)"
$ws.Range("L60").WrapText = $true

# --- Row heights (Excel recalculates these as best-fit heights on save) ---
$ws.Rows.Item(53).RowHeight = 17
$ws.Rows.Item(54).RowHeight = 17
$ws.Rows.Item(55).RowHeight = 17
$ws.Rows.Item(56).RowHeight = 17
$ws.Rows.Item(57).RowHeight = 17
$ws.Rows.Item(59).RowHeight = 17
$ws.Rows.Item(60).RowHeight = 34

# --- View state: scroll position / zoom / active cell selection ---
$ws.Application.ActiveWindow.Zoom = 107
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("L60").Select()

Write-Output "applied"
